$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source workbook stores every Price (D) / Volume(1h) (E) cell as plain
# text (inline strings), even when the Price text happens to look like a
# number (e.g. "216.61"). Most Volume cells are naturally safe because they
# carry padding spaces and a "%" sign, so Excel will not reinterpret them as
# numbers. A handful of the new Price values, though, are plain decimal-
# looking strings that Excel WOULD auto-convert to a Number on assignment
# (e.g. "6.90" -> 6.9, silently dropping the trailing zero). Force those
# specific cells to the Text number format first so the literal string is
# preserved exactly, matching the source data.
$textFormatCells = @("D5", "D6", "D8", "D10", "D15", "D17", "D20", "D21", "D23", "D27", "D28", "D33", "D38", "D42", "D43", "D46", "D49", "D51")
foreach ($addr in $textFormatCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "26.072.92"
$ws.Range("E2").Value = "  +0.54%  "
$ws.Range("D3").Value = "1.646.09"
$ws.Range("E3").Value = "  +0.80%  "
$ws.Range("E4").Value = "  +0.41%  "
$ws.Range("D5").Value = "216.61"
$ws.Range("E5").Value = "  +0.90%  "
$ws.Range("D6").Value = "0.508"
$ws.Range("E6").Value = "  +0.87%  "
$ws.Range("E7").Value = "  +0.38%  "
$ws.Range("D8").Value = "0.256"
$ws.Range("E8").Value = "  +0.36%  "
$ws.Range("E9").Value = "  +1.29%  "
$ws.Range("D10").Value = "19.67"
$ws.Range("E10").Value = "  +0.18%  "
$ws.Range("E11").Value = "  +0.58%  "
$ws.Range("D12").Value = "1.874.79"
$ws.Range("E12").Value = "  +0.80%  "
$ws.Range("D14").Value = "1.672.90"
$ws.Range("E14").Value = "  +2.29%  "
$ws.Range("D15").Value = "0.544"
$ws.Range("D16").Value = "0.0₃0764"
$ws.Range("E16").Value = "  +1.19%  "
$ws.Range("D17").Value = "63.44"
$ws.Range("E17").Value = "  +0.94%  "
$ws.Range("D18").Value = "26.083.67"
$ws.Range("E18").Value = "  +0.59%  "
$ws.Range("E19").Value = "  +0.46%  "
$ws.Range("D20").Value = "195.73"
$ws.Range("E20").Value = "  +1.42%  "
$ws.Range("D21").Value = "4.35"
$ws.Range("E21").Value = "  -0.58%  "
$ws.Range("E22").Value = "  -0.15%  "
$ws.Range("D23").Value = "6.23"
$ws.Range("E23").Value = "  -0.29%  "
$ws.Range("E24").Value = "  +0.58%  "
$ws.Range("E25").Value = "  +4.16%  "
$ws.Range("E26").Value = "  +0.81%  "
$ws.Range("D27").Value = "143.89"
$ws.Range("E27").Value = "  +0.69%  "
$ws.Range("D28").Value = "6.90"
$ws.Range("E28").Value = "  +0.76%  "
$ws.Range("E29").Value = "  +0.54%  "
$ws.Range("E30").Value = "  +0.96%  "
$ws.Range("E31").Value = "  -0.03%  "
$ws.Range("E32").Value = "  +1.94%  "
$ws.Range("D33").Value = "3.29"
$ws.Range("E33").Value = "  -0.43%  "
$ws.Range("E34").Value = "  -2.52%  "
$ws.Range("E36").Value = "  +0.60%  "
$ws.Range("D37").Value = "1.133.73"
$ws.Range("E37").Value = "  -0.29%  "
$ws.Range("D38").Value = "0.543"
$ws.Range("E38").Value = "  -1.28%  "
$ws.Range("E39").Value = "  -0.63%  "
$ws.Range("E40").Value = "  +0.30%  "
$ws.Range("D42").Value = "99.43"
$ws.Range("E42").Value = "  +0.27%  "
$ws.Range("D43").Value = "0.798"
$ws.Range("E43").Value = "  -0.76%  "
$ws.Range("D44").Value = "1.784.47"
$ws.Range("E44").Value = "  +0.83%  "
$ws.Range("D45").Value = "0.0₆0117"
$ws.Range("E45").Value = "  +3.79%  "
$ws.Range("D46").Value = "56.75"
$ws.Range("E46").Value = "  +0.75%  "
$ws.Range("E47").Value = "  +0.43%  "
$ws.Range("E48").Value = "  +0.41%  "
$ws.Range("D49").Value = "7.75"
$ws.Range("E49").Value = "  +1.84%  "
$ws.Range("E50").Value = "  +0.14%  "
$ws.Range("D51").Value = "0.0961"
$ws.Range("E51").Value = "  -0.07%  "
